$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, pushing existing data down
$ws.Rows.Item(1).Insert()

# New header row: "File Needed" / "Location" (bold + centered)
$ws.Range("A1").Value = "File Needed"
$ws.Range("B1").Value = "Location"

$a1 = $ws.Cells.Item(1, 1)
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108

$b1 = $ws.Cells.Item(1, 2)
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108

# Update the path values (now on rows 2 and 3) to the new Windows-style paths
$pairingsPath = "C:\Users\Jjacobson\Documents\Dropbox\Work-Personal Sync\Falsification\Scripts\Sourced\facial recognition\Facial-Validation\demo_faces\demo_filename_pairings.xlsx"
$imagesPath = "C:\Users\Jjacobson\Documents\Dropbox\Work-Personal Sync\Falsification\Scripts\Sourced\facial recognition\Facial-Validation\demo_faces\"

$ws.Range("B2").Value = $pairingsPath
$ws.Range("B3").Value = $imagesPath

# Extra reference to the pairings path further down the sheet
$ws.Range("B6").Value = $pairingsPath

# Column sizing to fit the new, much longer path text in column B
$ws.Columns.Item(2).AutoFit()

# Move the active selection
$ws.Range("B4").Select() | Out-Null

# Match the printed page orientation
$ws.PageSetup.Orientation = 1
